$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C:F for rows 2, 4-12 (row 3 unchanged) to reflect the re-ordered
# runs/balls/fours/sixes values described by the diff.
$data = @{
    2  = @("36", "16", "3", "3")
    4  = @("23", "16", "2", "1")
    5  = @("4", "10", "1", "0")
    6  = @("1", "5", "0", "0")
    7  = @("81", "51", "8", "3")
    8  = @("3", "2", "0", "0")
    9  = @("7", "9", "1", "0")
    10 = @("16", "22", "1", "0")
    11 = @("39", "34", "4", "2")
    12 = @("13", "12", "1", "0")
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $rng = $ws.Range("C$row`:F$row")
    $rng.NumberFormat = "@"
    for ($i = 0; $i -lt 4; $i++) {
        $ws.Cells.Item($row, 3 + $i).Value = $values[$i]
    }
}
